$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/000222a3b3e70f93181ffa723c3cbd634e306f77/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e035781949d1e586b4b193c798e16f8231ce7474/e2e/b.md."

# --- Overview sheet: row 3 (b.md) is now ready for handoff ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 16:44:47"

# --- zh-cn sheet: row 3 (b.md) handoff details ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# "False" must stay a text string (not Excel's boolean) - copy it from a cell
# that already holds the literal string "False" to avoid auto-coercion.
$zhcn.Range("O3").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-05 16:44:43"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.17

# --- de-de sheet: row 3 (b.md) handoff details ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-05 16:44:47"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.17
